$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(44326, 44327, 44328, 44329)
$startRow = 252

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i

    # Copy formatting (number format, font, borders, alignment) from the
    # cell directly above, then set the actual date-serial value.
    $ws.Cells.Item($r - 1, 1).Copy($ws.Cells.Item($r, 1))
    $ws.Cells.Item($r, 1).Value = $dates[$i]

    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
}
